# Auto-applied refresh of computed market-profit columns (H:N) across all crafting-job sheets.
# Values were refreshed by the scheduled Moogle/Universalis price-sync runner.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 1093.9333
$ws.Range("J33").Value = 2083.4
$ws.Range("L33").Value = 2083.4
$ws.Range("N33").Value = -2541.4
$ws.Range("H43").Value = 2723.0967
$ws.Range("I43").Value = 2505.5908
$ws.Range("J43").Value = 3254.7778
$ws.Range("K43").Value = 2505.5908
$ws.Range("L43").Value = 3254.7778
$ws.Range("M43").Value = -2436.5908
$ws.Range("N43").Value = -3392.7778
$ws.Range("H107").Value = 643.15
$ws.Range("I107").Value = 671
$ws.Range("J107").Value = 392.5
$ws.Range("K107").Value = 671
$ws.Range("L107").Value = 392.5
$ws.Range("M107").Value = 1249
$ws.Range("N107").Value = -4232.5
$ws.Range("H132").Value = 2846.3333
$ws.Range("I132").Value = 2806.5518
$ws.Range("K132").Value = 8419.6554
$ws.Range("M132").Value = -5889.6554
$ws.Range("H137").Value = 3024.093
$ws.Range("I137").Value = 2575.3447
$ws.Range("J137").Value = 3953.6428
$ws.Range("K137").Value = 7726.034100000001
$ws.Range("L137").Value = 11860.9284
$ws.Range("M137").Value = -5176.034100000001
$ws.Range("N137").Value = -16960.9284
$ws.Range("H138").Value = 4251.578
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4251.578
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 12754.734
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -23034.734

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 21600.63
$ws.Range("I32").Value = 15117
$ws.Range("K32").Value = 15117
$ws.Range("M32").Value = -14830
$ws.Range("H45").Value = 4115.4165
$ws.Range("I45").Value = 3990.6667
$ws.Range("J45").Value = 4157
$ws.Range("K45").Value = 3990.6667
$ws.Range("L45").Value = 4157
$ws.Range("M45").Value = -3613.6667
$ws.Range("N45").Value = -4911
$ws.Range("H61").Value = 7763.25
$ws.Range("I61").Value = 7431.282
$ws.Range("K61").Value = 7431.282
$ws.Range("M61").Value = -7219.282
$ws.Range("H74").Value = 5620.2256
$ws.Range("I74").Value = 3030.4348
$ws.Range("K74").Value = 3030.4348
$ws.Range("M74").Value = -2156.4348
$ws.Range("H77").Value = 5620.2256
$ws.Range("I77").Value = 3030.4348
$ws.Range("K77").Value = 15152.174
$ws.Range("M77").Value = -10784.174
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622
$ws.Range("H132").Value = 10989.454
$ws.Range("I132").Value = 3400
$ws.Range("J132").Value = 20096.8
$ws.Range("K132").Value = 10200
$ws.Range("L132").Value = 60290.39999999999
$ws.Range("M132").Value = -7670
$ws.Range("N132").Value = -65350.39999999999
$ws.Range("H136").Value = 7763.25
$ws.Range("I136").Value = 7431.282
$ws.Range("K136").Value = 22293.846
$ws.Range("M136").Value = -19743.846

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 2219.3635
$ws.Range("I86").Value = 2141.3
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 2141.3
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1018.3
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 2219.3635
$ws.Range("I89").Value = 2141.3
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 10706.5
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -5090.5
$ws.Range("N89").Value = -26232
$ws.Range("H99").Value = 3442.5715
$ws.Range("I99").Value = 1706.6666
$ws.Range("K99").Value = 1706.6666
$ws.Range("M99").Value = -208.6666
$ws.Range("H105").Value = 15494
$ws.Range("I105").Value = 12276
$ws.Range("J105").Value = 19999.2
$ws.Range("K105").Value = 12276
$ws.Range("L105").Value = 19999.2
$ws.Range("M105").Value = -10529
$ws.Range("N105").Value = -23493.2
$ws.Range("H107").Value = 2152.4614
$ws.Range("I107").Value = 2178.2
$ws.Range("J107").Value = 2066.6667
$ws.Range("K107").Value = 2178.2
$ws.Range("L107").Value = 2066.6667
$ws.Range("M107").Value = -258.1999999999998
$ws.Range("N107").Value = -5906.6667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("H18").Value = 65980
$ws.Range("J18").Value = 65980
$ws.Range("L18").Value = 65980
$ws.Range("N18").Value = -66440
$ws.Range("H31").Value = 8627.184999999999
$ws.Range("I31").Value = 4450.5
$ws.Range("J31").Value = 14370.125
$ws.Range("K31").Value = 4450.5
$ws.Range("L31").Value = 14370.125
$ws.Range("M31").Value = -4155.5
$ws.Range("N31").Value = -14960.125
$ws.Range("H34").Value = 8627.184999999999
$ws.Range("I34").Value = 4450.5
$ws.Range("J34").Value = 14370.125
$ws.Range("K34").Value = 4450.5
$ws.Range("L34").Value = 14370.125
$ws.Range("M34").Value = -4248.5
$ws.Range("N34").Value = -14774.125
$ws.Range("H94").Value = 1324.3914
$ws.Range("J94").Value = 1409.2632
$ws.Range("L94").Value = 1409.2632
$ws.Range("N94").Value = -2311.2632
$ws.Range("H105").Value = 2146.4546
$ws.Range("I105").Value = 2261.1
$ws.Range("K105").Value = 2261.1
$ws.Range("M105").Value = -514.0999999999999
$ws.Range("H107").Value = 1538
$ws.Range("I107").Value = 1291.8636
$ws.Range("J107").Value = 2139.6667
$ws.Range("K107").Value = 1291.8636
$ws.Range("L107").Value = 2139.6667
$ws.Range("M107").Value = 628.1364000000001
$ws.Range("N107").Value = -5979.6667
$ws.Range("H132").Value = 6402.4
$ws.Range("I132").Value = 5093.6816
$ws.Range("K132").Value = 15281.0448
$ws.Range("M132").Value = -12751.0448

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item(5)
$ws.Range("H32").Value = 4473.9565
$ws.Range("I32").Value = 1750
$ws.Range("J32").Value = 4733.381
$ws.Range("K32").Value = 5250
$ws.Range("L32").Value = 14200.143
$ws.Range("M32").Value = -4967
$ws.Range("N32").Value = -14766.143
$ws.Range("H98").Value = 390.92856
$ws.Range("I98").Value = 344.66666
$ws.Range("J98").Value = 425.625
$ws.Range("K98").Value = 1033.99998
$ws.Range("L98").Value = 1276.875
$ws.Range("M98").Value = 464.0000199999999
$ws.Range("N98").Value = -4272.875
$ws.Range("H138").Value = 4134.75
$ws.Range("J138").Value = 4364.6665
$ws.Range("L138").Value = 13093.9995
$ws.Range("N138").Value = -23373.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("H63").Value = 45000
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 45000
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 7823.35
$ws.Range("I122").Value = 4582.2
$ws.Range("J122").Value = 11064.5
$ws.Range("K122").Value = 13746.6
$ws.Range("L122").Value = 33193.5
$ws.Range("M122").Value = -11296.6
$ws.Range("N122").Value = -38093.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 336334.66
$ws.Range("I7").Value = 502002
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 502002
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -501890
$ws.Range("N7").Value = -5224
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 3420.5
$ws.Range("I46").Value = 3389.4285
$ws.Range("K46").Value = 3389.4285
$ws.Range("M46").Value = -3201.4285
$ws.Range("H55").Value = 958.1111
$ws.Range("I55").Value = 181.5
$ws.Range("J55").Value = 1579.4
$ws.Range("K55").Value = 181.5
$ws.Range("L55").Value = 1579.4
$ws.Range("M55").Value = -8.5
$ws.Range("N55").Value = -1925.4
$ws.Range("H93").Value = 3450.8572
$ws.Range("I93").Value = 831.2
$ws.Range("K93").Value = 831.2
$ws.Range("M93").Value = 416.8
$ws.Range("H126").Value = 336334.66
$ws.Range("I126").Value = 502002
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 1506006
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -1503536
$ws.Range("N126").Value = -19940

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("H44").Value = 143666.67
$ws.Range("J44").Value = 143666.67
$ws.Range("L44").Value = 143666.67
$ws.Range("N44").Value = -144774.67
$ws.Range("H100").Value = 1455.7
$ws.Range("I100").Value = 910.8
$ws.Range("J100").Value = 2000.6
$ws.Range("K100").Value = 1821.6
$ws.Range("L100").Value = 4001.2
$ws.Range("M100").Value = -1280.6
$ws.Range("N100").Value = -5083.2
$ws.Range("H107").Value = 3901.75
$ws.Range("I107").Value = 3306.8147
$ws.Range("K107").Value = 9920.444100000001
$ws.Range("M107").Value = -8000.444100000001
$ws.Range("H126").Value = 3608.9048
$ws.Range("I126").Value = 3005.7646
$ws.Range("J126").Value = 6172.25
$ws.Range("K126").Value = 9017.293799999999
$ws.Range("L126").Value = 18516.75
$ws.Range("M126").Value = -6547.293799999999
$ws.Range("N126").Value = -23456.75
$ws.Range("H132").Value = 2235.9824
$ws.Range("I132").Value = 1989.1837
$ws.Range("K132").Value = 5967.551100000001
$ws.Range("M132").Value = -3437.551100000001
$ws.Range("H136").Value = 3261.639
$ws.Range("I136").Value = 2877.742
$ws.Range("J136").Value = 5641.8
$ws.Range("K136").Value = 8633.226000000001
$ws.Range("L136").Value = 16925.4
$ws.Range("M136").Value = -6083.226000000001
$ws.Range("N136").Value = -22025.4
